$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '70.203.82'
$ws.Range('E2').Value = '  +0.74%  '
$ws.Range('D3').Value = '3.503.71'
$ws.Range('E3').Value = '  +0.10%  '
$ws.Range('E4').Value = '  -0.07%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '604.82'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +0.11%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '172.94'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +1.04%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.611'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  -0.64%  '
$ws.Range('D8').Value = '3.498.00'
$ws.Range('E8').Value = '  +0.06%  '
$ws.Range('E9').Value = '  +0.06%  '
$ws.Range('E10').Value = '  -2.36%  '
$ws.Range('E11').Value = '  +7.03%  '
$ws.Range('E12').Value = '  +1.06%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '46.05'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -2.01%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.0000277'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -0.91%  '
$ws.Range('D15').Value = '4.073.00'
$ws.Range('E15').Value = '  +0.06%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '615.01'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -0.85%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '8.28'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -1.14%  '
$ws.Range('D18').Value = '3.500.94'
$ws.Range('E18').Value = '  +0.05%  '
$ws.Range('D19').Value = '70.160.56'
$ws.Range('E19').Value = '  +0.66%  '
$ws.Range('E20').Value = '  +1.09%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '17.51'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +1.48%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.879'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -0.22%  '
$ws.Range('E23').Value = '  -8.00%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '99.21'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +3.17%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '15.59'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -1.18%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '3.72'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -3.09%  '
$ws.Range('E27').Value = '  -0.04%  '
$ws.Range('E28').Value = '  -0.97%  '
$ws.Range('E29').Value = '  +2.87%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '9.02'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -1.36%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '658.96'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +17.29%  '
$ws.Range('B32').Value = 'Filecoin'
$ws.Range('C32').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '8.04'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -3.98%  '
$ws.Range('B33').Value = 'Stacks'
$ws.Range('C33').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '2.96'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -2.77%  '
$ws.Range('E34').Value = '  -4.23%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '6.80'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -1.86%  '
$ws.Range('E36').Value = '  -1.36%  '
$ws.Range('E37').Value = '  -0.17%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '3.51'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +0.90%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.0476'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +6.37%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '56.62'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -0.62%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '1.00'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -0.04%  '
$ws.Range('E42').Value = '  +1.73%  '
$ws.Range('D43').Value = '3.352.34'
$ws.Range('E43').Value = '  +0.86%  '
$ws.Range('D44').Value = '0.0₃0736'
$ws.Range('E44').Value = '  +4.10%  '
$ws.Range('E45').Value = '  -4.71%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '2.90'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -2.38%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '31.87'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -3.03%  '
$ws.Range('E48').Value = '  -2.30%  '
$ws.Range('E49').Value = '  +0.99%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '132.82'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -1.50%  '
